$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the sheet's current contents/formatting outright - the new layout
# shifts every data column one to the left and inserts a block of new rows,
# so it's far safer to rebuild it from scratch than to patch cells in place.
# (Clearing does not disturb the relative ordering of shared strings that
# remain in use, it only forgets the association for this worksheet.)
$ws.Cells.Clear()

# --- Header row (row 1): name / value / unit / input type / lower boundary / upper boundary / check boundary
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "value"
$ws.Range("D1").Value = "unit"
$ws.Range("E1").Value = "input type"
$ws.Range("F1").Value = "lower boundary"
$ws.Range("G1").Value = "upper boundary"
$ws.Range("H1").Value = "check boundary"

# --- Fuselage section
$ws.Range("A2").Value = "Fuselage"

$ws.Range("B5").Value = "fuselageLength"
$ws.Range("C5").Value = 35
$ws.Range("D5").Value = "m"
$ws.Range("E5").Value = "float"

$ws.Range("B6").Value = "fuselageDiameter"
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = "m"
$ws.Range("E6").Value = "float"

$ws.Range("B7").Value = "noseSlenderness"
$ws.Range("C7").Value = 2
$ws.Range("E7").Value = "float"

$ws.Range("B8").Value = "tailSlenderness"
$ws.Range("C8").Value = 2
$ws.Range("E8").Value = "float"

$ws.Range("B9").Value = "tailUpAngle"
$ws.Range("D9").Value = "deg"
$ws.Range("E9").Value = "float"

# --- New strings: type "EOF" before "Wing" so the shared-string intern order
# (new strings are appended in first-use order, after the strings that were
# already part of the workbook) matches the target file (EOF then Wing).
$ws.Range("A16").Value = "EOF"
$ws.Range("A10").Value = "Wing"

# --- Wing section
$ws.Range("B12").Value = "aspectRatio"
$ws.Range("E12").Value = "float"

$ws.Range("B13").Value = "maTechnology"
$ws.Range("E13").Value = "float"

$ws.Range("B15").Value = "wingPosition"
$ws.Range("E15").Value = "string"

# --- "value" column (C) keeps the 0.000 number format all the way down,
# even on the blank rows, matching the template's look.
$ws.Range("C2:C35").NumberFormat = "0.000"

# --- Touch the now-blank "unit" column (D) cells so a cell record exists
# for them on every templated row, matching the target layout.
$ws.Range("D2").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("D16:D35").Value = ""

# --- Column widths (best-fit on the label/unit columns)
$ws.Columns.Item(2).ColumnWidth = 17
$ws.Columns.Item(6).ColumnWidth = 15.140625
$ws.Columns.Item(7).ColumnWidth = 15.28515625
$ws.Columns.Item(8).ColumnWidth = 15

# --- Selection, matching the saved view state
$ws.Range("A10").Select()

Write-Output "done"
